$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.895.45"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "1.702.68"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("E4").Value = "  -0.60%  "

$ws.Range("D5").Value = "'315.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("E6").Value = "  -0.65%  "

$ws.Range("D7").Value = "'0.4046"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.48%  "

$ws.Range("D8").Value = "'0.4058"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.55%  "

$ws.Range("D9").Value = "'1.003"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.69%  "

$ws.Range("D10").Value = "'53.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.22%  "

$ws.Range("D11").Value = "'1.468"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.47%  "

$ws.Range("D12").Value = "'0.08809"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("D13").Value = "'25.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.81%  "

$ws.Range("D14").Value = "'7.526"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.49%  "

$ws.Range("D15").Value = "'8.047"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.38%  "

$ws.Range("D16").Value = "'0.00001351"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.54%  "

$ws.Range("D17").Value = "1.726.42"
$ws.Range("E17").Value = "  +1.53%  "

$ws.Range("D18").Value = "'96.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.18%  "

$ws.Range("D19").Value = "'0.07169"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.63%  "

$ws.Range("D20").Value = "'21.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.16%  "

$ws.Range("D21").Value = "'7.234"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.55%  "

$ws.Range("E22").Value = "  -0.84%  "

$ws.Range("E23").Value = "  +1.23%  "

$ws.Range("D24").Value = "24.897.24"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("E25").Value = "  -1.18%  "

$ws.Range("D26").Value = "'6.835"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +30.66%  "

$ws.Range("D27").Value = "'2.889"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.69%  "

$ws.Range("E28").Value = "  +0.82%  "

$ws.Range("D29").Value = "'164.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("D30").Value = "'145.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.44%  "

$ws.Range("D31").Value = "'8.238"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.87%  "

$ws.Range("B32").Value = "WEMIXTOKEN"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "'2.270"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.98%  "

$ws.Range("B33").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C33").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D33").Value = "1.917.29"
$ws.Range("E33").Value = "  +1.44%  "

$ws.Range("D34").Value = "'0.08805"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.65%  "

$ws.Range("D35").Value = "'0.03201"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.54%  "

$ws.Range("D36").Value = "'7.312"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.44%  "

$ws.Range("D37").Value = "'1.014"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.07%  "

$ws.Range("D38").Value = "'0.2844"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.14%  "

$ws.Range("D39").Value = "'0.8450"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.65%  "

$ws.Range("D40").Value = "'10.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.22%  "

$ws.Range("D41").Value = "'0.09399"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.30%  "

$ws.Range("D42").Value = "'14.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.61%  "

$ws.Range("D43").Value = "'18.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.33%  "

$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("D45").Value = "'2.714"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.43%  "

$ws.Range("D46").Value = "'0.7429"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.87%  "

$ws.Range("D47").Value = "'4.241"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.41%  "

$ws.Range("E48").Value = "  +4.11%  "

$ws.Range("D49").Value = "'1.003"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").Value = "'142.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.40%  "

$ws.Range("E51").Value = "  +3.61%  "
